$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '65.677.33'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.651.32'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '597.88'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '156.61'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.627'
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '0.0000198'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '28.66'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').Value = '3.127.87'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '65.502.76'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '2.670.92'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '12.61'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = '348.72'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '68.96'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').Value = '0.0000114'
$ws.Range('E24').Value = '  +2.99%  '
$ws.Range('D25').Value = '9.65'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('E26').Value = '  +3.88%  '
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '531.92'
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.12'
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').Value = '6.42'
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').Value = '5.41'
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').Value = '20.37'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = '155.88'
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '161.30'
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '0.0607'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').Value = '2.28'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = '22.61'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').Value = '0.635'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('D49').Value = '0.0995'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('E50').Value = '  +7.70%  '
$ws.Range('D51').Value = '19.74'
$ws.Range('E51').Value = '  -2.66%  '
